# Tabela de atores e respetivos casos de uso - revisões de diagrama e estrutura do relatorio
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows belonging to the old "questionario" entries that were merged/renamed
# ("Eliminar questionario" and "Alterar questionario"), shifting everything below up by two rows.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(8).Delete()

# Fix/normalize the remaining case-use descriptions (capitalisation & wording updates)
$ws.Range("B7").Value = "Responder Questionario"
$ws.Range("B8").Value = "Consultar Questionario"
$ws.Range("B9").Value = "Criar Turista"
$ws.Range("B10").Value = "Consultar Turista"
$ws.Range("B11").Value = "Alterar Turista"
$ws.Range("B12").Value = "Eliminar Turista"
$ws.Range("B13").Value = "Visualizar Trilho Recomendado"
$ws.Range("B15").Value = "Consultar Estabelecimentos para Descansar"
$ws.Range("B16").Value = "Visualizar Condições Metereologicas para o dia selecionado"
$ws.Range("B17").Value = "Selecionar o dia que deseja fazer o trilho"
$ws.Range("B20").Value = "Selecionar Guia"
$ws.Range("B21").Value = "Alterar Guia"
$ws.Range("B22").Value = "Escolher Guia"
$ws.Range("B23").Value = "Mostrar Contactos"
$ws.Range("B24").Value = "Consultar Testes Trilho"

# Update the view's active selection to match the saved state
$ws.Range("B25").Select()
